$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NCT(2.260085008522421, 1.4211906422018523, -0.5172128520750621, 2.214272917878322)"
$ws.Range("C2").Value = "NCT(3.2260123264523193, 1.1015571026262847, 3.4806849659485977, 4.5238235005647365)"
$ws.Range("D2").Value = "JSU(-0.9198199206443796, 1.2441341260418963, 0.6722253377479606, 3.11334946646099)"
$ws.Range("E2").Value = "JSU(-1.223026086395933, 1.0028476383675553, 3.0610907124776148, 3.3015988855807317)"
